# Swap the order of names in the "Recorded By" (column G) cells that list
# both "dnasr281@gmail.com" and "System", changing
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
# This mirrors the change seen in the XML diff for the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
